# Update Leve market-price figures (currentAveragePrice / NQ / HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ)
# across the Unicorn_Profits sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Sheets("ALC")
$ws.Range("H76").Value = 7429.3438
$ws.Range("I76").Value = 9573.823
$ws.Range("J76").Value = 4998.933
$ws.Range("K76").Value = 9573.823
$ws.Range("L76").Value = 4998.933
$ws.Range("M76").Value = -9258.823
$ws.Range("N76").Value = -5628.933
$ws.Range("H79").Value = 7429.3438
$ws.Range("I79").Value = 9573.823
$ws.Range("J79").Value = 4998.933
$ws.Range("K79").Value = 9573.823
$ws.Range("L79").Value = 4998.933
$ws.Range("M79").Value = -8481.823
$ws.Range("N79").Value = -7182.933
$ws.Range("H132").Value = 2793.8833
$ws.Range("I132").Value = 1548.7693
$ws.Range("J132").Value = 5106.2383
$ws.Range("K132").Value = 4646.3079
$ws.Range("L132").Value = 15318.7149
$ws.Range("M132").Value = -2116.3079
$ws.Range("N132").Value = -20378.7149
$ws.Range("H141").Value = 459.45312
$ws.Range("I141").Value = 455
$ws.Range("K141").Value = 1365
$ws.Range("M141").Value = 3815

# --- Sheet: ARM ---
$ws = $wb.Sheets("ARM")
$ws.Range("H32").Value = 5928.344
$ws.Range("I32").Value = 3883.7036
$ws.Range("K32").Value = 3883.7036
$ws.Range("M32").Value = -3596.7036
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H45").Value = 1040.8334
$ws.Range("I45").Value = 961.4
$ws.Range("J45").Value = 1438
$ws.Range("K45").Value = 961.4
$ws.Range("L45").Value = 1438
$ws.Range("M45").Value = -584.4
$ws.Range("N45").Value = -2192
$ws.Range("H61").Value = 282010.38
$ws.Range("I61").Value = 197681.02
$ws.Range("J61").Value = 530560.1
$ws.Range("K61").Value = 197681.02
$ws.Range("L61").Value = 530560.1
$ws.Range("M61").Value = -197469.02
$ws.Range("N61").Value = -530984.1
$ws.Range("H97").Value = 783.73914
$ws.Range("I97").Value = 858.2778
$ws.Range("J97").Value = 515.4
$ws.Range("K97").Value = 858.2778
$ws.Range("L97").Value = 515.4
$ws.Range("M97").Value = -362.2778
$ws.Range("N97").Value = -1507.4
$ws.Range("H136").Value = 282010.38
$ws.Range("I136").Value = 197681.02
$ws.Range("J136").Value = 530560.1
$ws.Range("K136").Value = 593043.0599999999
$ws.Range("L136").Value = 1591680.3
$ws.Range("M136").Value = -590493.0599999999
$ws.Range("N136").Value = -1596780.3

# --- Sheet: BSM ---
$ws = $wb.Sheets("BSM")
$ws.Range("H105").Value = 1865.3636
$ws.Range("I105").Value = 1389.3846
$ws.Range("K105").Value = 1389.3846
$ws.Range("M105").Value = 357.6153999999999

# --- Sheet: CRP ---
$ws = $wb.Sheets("CRP")
$ws.Range("H31").Value = 2777.6365
$ws.Range("I31").Value = 1849.2106
$ws.Range("J31").Value = 4852.9414
$ws.Range("K31").Value = 1849.2106
$ws.Range("L31").Value = 4852.9414
$ws.Range("M31").Value = -1554.2106
$ws.Range("N31").Value = -5442.9414
$ws.Range("H34").Value = 2777.6365
$ws.Range("I34").Value = 1849.2106
$ws.Range("J34").Value = 4852.9414
$ws.Range("K34").Value = 1849.2106
$ws.Range("L34").Value = 4852.9414
$ws.Range("M34").Value = -1647.2106
$ws.Range("N34").Value = -5256.9414
$ws.Range("H58").Value = 2688.9219
$ws.Range("I58").Value = 2681.9019
$ws.Range("J58").Value = 2716.4614
$ws.Range("K58").Value = 2681.9019
$ws.Range("L58").Value = 2716.4614
$ws.Range("M58").Value = -2478.9019
$ws.Range("N58").Value = -3122.4614
$ws.Range("H136").Value = 2688.9219
$ws.Range("I136").Value = 2681.9019
$ws.Range("J136").Value = 2716.4614
$ws.Range("K136").Value = 8045.7057
$ws.Range("L136").Value = 8149.3842
$ws.Range("M136").Value = -5495.7057
$ws.Range("N136").Value = -13249.3842

# --- Sheet: CUL ---
$ws = $wb.Sheets("CUL")
$ws.Range("H81").Value = 5660
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 5660
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 16980
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -19226
$ws.Range("H84").Value = 5660
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 5660
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 50940
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -62172
$ws.Range("H113").Value = 572.8461
$ws.Range("I113").Value = 567.875
$ws.Range("J113").Value = 580.8
$ws.Range("K113").Value = 1703.625
$ws.Range("L113").Value = 1742.4
$ws.Range("M113").Value = 466.375
$ws.Range("N113").Value = -6082.4
$ws.Range("H131").Value = 1157.1562
$ws.Range("I131").Value = 1028
$ws.Range("J131").Value = 1200.2084
$ws.Range("K131").Value = 3084
$ws.Range("L131").Value = 3600.6252
$ws.Range("M131").Value = 1956
$ws.Range("N131").Value = -13680.6252

# --- Sheet: GSM ---
$ws = $wb.Sheets("GSM")
$ws.Range("H97").Value = 785.1852
$ws.Range("I97").Value = 617.3913
$ws.Range("J97").Value = 1750
$ws.Range("K97").Value = 617.3913
$ws.Range("L97").Value = 1750
$ws.Range("M97").Value = -121.3913
$ws.Range("N97").Value = -2742

# --- Sheet: LTW ---
$ws = $wb.Sheets("LTW")
$ws.Range("H68").Value = 3154.6667
$ws.Range("I68").Value = 2644
$ws.Range("J68").Value = 3410
$ws.Range("K68").Value = 2644
$ws.Range("L68").Value = 3410
$ws.Range("M68").Value = -1895
$ws.Range("N68").Value = -4908
$ws.Range("H71").Value = 3154.6667
$ws.Range("I71").Value = 2644
$ws.Range("J71").Value = 3410
$ws.Range("K71").Value = 13220
$ws.Range("L71").Value = 17050
$ws.Range("M71").Value = -9476
$ws.Range("N71").Value = -24538
$ws.Range("H100").Value = 83339280
$ws.Range("I100").Value = 11399.8
$ws.Range("J100").Value = 142859200
$ws.Range("K100").Value = 11399.8
$ws.Range("L100").Value = 142859200
$ws.Range("M100").Value = -10858.8
$ws.Range("N100").Value = -142860282
$ws.Range("H122").Value = 2685.7144
$ws.Range("I122").Value = 2825
$ws.Range("K122").Value = 8475
$ws.Range("M122").Value = -6025

# --- Sheet: WVR ---
$ws = $wb.Sheets("WVR")
$ws.Range("H107").Value = 300
$ws.Range("I107").Value = 300
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -4740
$ws.Range("H113").Value = 417.54544
$ws.Range("I113").Value = 405.375
$ws.Range("K113").Value = 1216.125
$ws.Range("M113").Value = 953.875
$ws.Range("H126").Value = 2299.52
$ws.Range("I126").Value = 2304
$ws.Range("J126").Value = 2294.6667
$ws.Range("K126").Value = 6912
$ws.Range("L126").Value = 6884.000100000001
$ws.Range("M126").Value = -4442
$ws.Range("N126").Value = -11824.0001

